$d = $word.ActiveDocument

# 1) Expand " uit gezet?" into " uit gezet, ook slim zoeken content plugin?"
$null = $d.Content.Find.Execute(" uit gezet?", $true, $false, $false, $false, $false, $true, 1, $false, " uit gezet, ook slim zoeken content plugin?", 2)

# 2) Split the trailing "plugin?" off into its own runs ("plugin" and "?") while keeping
#    identical run formatting - force a transient font-size change then revert it so the
#    engine records the run boundary without leaving a formatting diff behind.
$scoped = $d.Range(1300, $d.Content.End)
$null = $scoped.Find.Execute("plugin", $true, $false, $false, $false, $false, $true)
$scoped.Font.Size = 11
$scoped.Font.Size = 10

# 3) Re-find the exact "plugin" run and wrap it with proofErr spellStart/spellEnd markers.
$scoped2 = $d.Range(1300, $d.Content.End)
$null = $scoped2.Find.Execute("plugin", $true, $false, $false, $false, $false, $true)
$pluginStart = $scoped2.Start
$pluginEnd = $scoped2.End

# 4) Move the existing _GoBack bookmark from the "Slim zoeken indexen..." paragraph to
#    right after the new "?" run.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$scoped3 = $d.Range(1300, $d.Content.End)
$null = $scoped3.Find.Execute("plugin?", $true, $false, $false, $false, $false, $true)
$afterQuestion = $d.Range($scoped3.End, $scoped3.End)
$d.Bookmarks.Add("_GoBack", $afterQuestion)
